$wb = $excel.ActiveWorkbook

# Commit: "Update chuc nang anh dai dien" / "them vao 2 bang HV va GV 1
# truong Image" -> add an "Image" column to both the GiangVien (GV) and
# HocVien (HV) tables. In both sheets the new column is inserted right
# before the existing last (department/class) column, pushing it one
# column to the right.

# --- 2GiangVien (GV) sheet ---
$wsGV = $wb.Worksheets.Item("2GiangVien")
$wsGV.Columns("H").Insert() | Out-Null
$wsGV.Columns("H").ColumnWidth = $wsGV.Columns("G").ColumnWidth
$wsGV.Range("H1:H1048576").Select() | Out-Null

# --- 7HocVien (HV) sheet ---
$wsHV = $wb.Worksheets.Item("7HocVien")
$wsHV.Columns("H").Insert() | Out-Null
$wsHV.Columns("H").ColumnWidth = $wsHV.Columns("G").ColumnWidth
$wsHV.Range("H7").Select() | Out-Null

$wsHV.Activate() | Out-Null
